# Update gh-pages output data (regenerated scrape) - apply the new
# "想去人数" (attendance) figures to the affected rows across sheets.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 7787
$ws1.Range("F24").Value = 76
$ws1.Range("F29").Value = 4079
$ws1.Range("F30").Value = 2
$ws1.Range("F35").Value = 363
$ws1.Range("F36").Value = 1413
$ws1.Range("F40").Value = 3237

# Sheet "本地生活" (sheet3)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 1355

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1355
$ws4.Range("F5").Value = 7787
$ws4.Range("F26").Value = 76
$ws4.Range("F30").Value = 4079
$ws4.Range("F34").Value = 363
$ws4.Range("F36").Value = 1413
$ws4.Range("F41").Value = 3237
